$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 392-393, shifting the existing rows 392-397 down to 394-399.
$ws.Rows("392:393").Insert()

# --- New row 392 ---
$ws.Range("A392").Value = 6
$ws.Range("B392").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C392").Value = "Metropolitana"
$ws.Range("D392").Value = 44595
$ws.Range("E392").Value = 13
$ws.Range("F392").Value = 100112030
$ws.Range("G392").Value = "Poroto granado"
$ws.Range("H392").Value = "Sin especificar"
$ws.Range("I392").Value = "Primera"
$ws.Range("J392").Value = 1000
$ws.Range("K392").Value = 20000
$ws.Range("L392").Value = 23000
$ws.Range("M392").Value = 21200
$ws.Range("N392").Value = "$/saco 25 kilos"
$ws.Range("O392").Value = "Región Metropolitana"
$ws.Range("P392").Value = 848
$ws.Range("Q392").Value = 25
$ws.Range("R392").Value = "Hortaliza"

# --- New row 393 ---
$ws.Range("A393").Value = 6
$ws.Range("B393").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C393").Value = "Metropolitana"
$ws.Range("D393").Value = 44595
$ws.Range("E393").Value = 13
$ws.Range("F393").Value = 100112030
$ws.Range("G393").Value = "Poroto granado"
$ws.Range("H393").Value = "Sin especificar"
$ws.Range("I393").Value = "Primera"
$ws.Range("J393").Value = 750
$ws.Range("K393").Value = 20000
$ws.Range("L393").Value = 25000
$ws.Range("M393").Value = 22667
$ws.Range("N393").Value = "$/saco 25 kilos"
$ws.Range("O393").Value = "Región del Maule"
$ws.Range("P393").Value = 907
$ws.Range("Q393").Value = 25
$ws.Range("R393").Value = "Hortaliza"
